$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (Session 7 / pre-tx phase parent outcome measures)
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 4

# Update the selected cell to reflect the new active selection after entry
$ws.Range("B3").Select()
